$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replicate the "spacer + Assignment + Grade + Comments" block (cols Y:AB,
# rows 1-12) six times into AC:AF, AG:AJ, AK:AN, AO:AR, AS:AV, AW:AZ so the new
# columns inherit the same styles (header fill/border, spacer fill, data
# border) as the existing assignment blocks. ---
$template = $ws.Range("Y1:AB12")
$template.Copy($ws.Range("AC1"))
$template.Copy($ws.Range("AG1"))
$template.Copy($ws.Range("AK1"))
$template.Copy($ws.Range("AO1"))
$template.Copy($ws.Range("AS1"))
$template.Copy($ws.Range("AW1"))

# Row 1 (headers) keeps the same Assignment / Grade / Comments labels in every
# new block - already correct after the copy above. Only the row height
# changes (the header row grows to fit two lines of wrapped text).
$ws.Range("A1:AZ1").RowHeight = 30

# --- Row 2: new per-assignment grade data -------------------------------
# Block "7 CPP"
$ws.Range("AD2").Value = "7 CPP"
$ws.Range("AE2").Value = 85
$ws.Range("AF2").Value = "Good!, Please read my comments"

# Block "8 CPP"
$ws.Range("AH2").Value = "8 CPP"
$ws.Range("AI2").Value = 100
$ws.Range("AJ2").Value = "Excellent!"

# Block "9 CPP"
$ws.Range("AL2").Value = "9 CPP"
$ws.Range("AM2").Value = 99
$ws.Range("AN2").Value = "Excellent! (watch out for those memory leaks)"

# Block "10 CPP"
$ws.Range("AP2").Value = "10 CPP"
$ws.Range("AQ2").Value = 91
$ws.Range("AR2").Value = "Very good!"

# Block "11 CPP"
$ws.Range("AT2").Value = "11 CPP"
$ws.Range("AU2").Value = 71
$ws.Range("AV2").Value = "See my Comments"

# Block "12 CPP"
$ws.Range("AX2").Value = "12 CPP"
$ws.Range("AY2").Value = 100
$ws.Range("AZ2").Value = "Excellent!"

# --- Selection / view matches where the user ended up after adding columns --
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("AC2:AZ2").Select()
